$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) from numeric auto-coercion while we write plain decimal-looking text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '44.034.95'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '2.360.89'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '74.13'
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.593'
$ws.Range("E9").Value = '  +7.25%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '32.16'
$ws.Range("E12").Value = '  +7.74%  '
$ws.Range("E13").Value = '  +7.06%  '
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").Value = '2.711.48'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '16.58'
$ws.Range("E16").Value = '  -1.97%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '2.347.37'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").Value = '43.941.44'
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("D20").Value = '7.03'
$ws.Range("E20").Value = '  +8.47%  '
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").Value = '77.32'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '258.65'
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").Value = '  +23.19%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").Value = '3.66'
$ws.Range("E26").Value = '  -4.55%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '2.49'
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").Value = '10.80'
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '22.80'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = '175.53'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("D35").Value = '5.63'
$ws.Range("E35").Value = '  +8.18%  '
$ws.Range("D36").Value = '5.23'
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '3.77'
$ws.Range("E37").Value = '  -3.86%  '
$ws.Range("D38").Value = '6.37'
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").Value = '2.36'
$ws.Range("E39").Value = '  -3.23%  '
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("E41").Value = '  +14.09%  '
$ws.Range("E42").Value = '  +10.25%  '
$ws.Range("D43").Value = '19.01'
$ws.Range("E43").Value = '  -5.16%  '
$ws.Range("D44").Value = '9.01'
$ws.Range("E44").Value = '  +1.77%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '4.76'
$ws.Range("E46").Value = '  +6.20%  '
$ws.Range("D47").Value = '59.34'
$ws.Range("E47").Value = '  +12.90%  '
$ws.Range("E48").Value = '  +6.19%  '
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").Value = '101.03'
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("E51").Value = '  +0.06%  '

# Restore original (default) cell style now that text values are safely stored
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Applied cryptos list update"
